# Add explanatory "comments" (column D) to the GME configuration parameters
# sheet (WiFi), describing each configuration field, and update the note
# about the config-file "0=Enabled 1=Disabled" cell with additional remarks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WiFi")

# The order below matters: it reproduces the exact order in which new
# shared-string entries were appended to xl/sharedStrings.xml in the
# target workbook.
$ws.Range("D11").Value = "ssid of the AP to which GME should connect"
$ws.Range("D13").Value = "psk of the AP to which GME should connect"
$ws.Range("D3").Value  = "don't know"
$ws.Range("D14").Value = "GME will acquire an IP address from a dhcp server when connected to the AP? 0:no, 1:yes "
$ws.Range("D15").Value = "Static IP address of the GME when dhcp is disabled"
$ws.Range("D17").Value = "Gateway of the GME when dhcp is disabled"
$ws.Range("D16").Value = "Netmask of the GME when dhcp is disabled"
$ws.Range("D2").Value  = "ssid of the GME"
$ws.Range("D4").Value  = "psk to access GME as AP"
$ws.Range("D21").Value = "DNS"
$ws.Range("D19").Value = "APN"
$ws.Range("D12").Value = "don't know (do we support different types of encryption?)"
$ws.Range("D9").Value  = "IP address leased by a dhcp server"
$ws.Range("D6").Value  = "IP address of the GME working as AP"
$ws.Range("D7").Value  = "netmask of the GME working as AP"

# D8 previously held the "0=Enabled  1=Disabled" note; extend it with an
# extra remark on its own line, wrap the text and grow the row to fit it.
$ws.Range("D8").Value = "secondo me non serve, perché il modulo non avrà un server dhcp a bordo per distribuire indirizzi" + [char]10 + "0=Enabled  1=Disabled"
$ws.Range("D8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 43.2
